$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF), written as text like the
# other header-row cells (B1:H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Reuse the exact same formatting as the rest of the header row (bold,
# bordered, centered) by copying the format from H1 onto I1:J1, rather than
# fabricating a brand-new style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data for columns I (I0) and J (IF), rows 2-5.
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 6
